$d = $word.ActiveDocument
$t = $d.Tables(1)
$vtab = [char]11

$cell = $t.Cell(1,1)
$cell.Range.Text = "57 x 86" + $vtab + "  8    6" + $vtab + "  ----" + $vtab + "5|    |" + $vtab + "7|    |"

$cell = $t.Cell(1,2)
$cell.Range.Text = "28 x 11" + $vtab + "  1    1" + $vtab + "  ----" + $vtab + "2|    |" + $vtab + "8|    |"

$cell = $t.Cell(1,3)
$cell.Range.Text = "66 x 92" + $vtab + "  9    2" + $vtab + "  ----" + $vtab + "6|    |" + $vtab + "6|    |"

$cell = $t.Cell(2,1)
$cell.Range.Text = "75 x 82" + $vtab + "  8    2" + $vtab + "  ----" + $vtab + "7|    |" + $vtab + "5|    |"

$cell = $t.Cell(2,2)
$cell.Range.Text = "32 x 91" + $vtab + "  9    1" + $vtab + "  ----" + $vtab + "3|    |" + $vtab + "2|    |"

$cell = $t.Cell(2,3)
$cell.Range.Text = "67 x 89" + $vtab + "  8    9" + $vtab + "  ----" + $vtab + "6|    |" + $vtab + "7|    |"

$cell = $t.Cell(3,1)
$cell.Range.Text = "67 x 93" + $vtab + "  9    3" + $vtab + "  ----" + $vtab + "6|    |" + $vtab + "7|    |"

$cell = $t.Cell(3,2)
$cell.Range.Text = "71 x 30" + $vtab + "  3    0" + $vtab + "  ----" + $vtab + "7|    |" + $vtab + "1|    |"

$cell = $t.Cell(3,3)
$cell.Range.Text = "47 x 58" + $vtab + "  5    8" + $vtab + "  ----" + $vtab + "4|    |" + $vtab + "7|    |"

$cell = $t.Cell(4,1)
$cell.Range.Text = "29 x 93" + $vtab + "  9    3" + $vtab + "  ----" + $vtab + "2|    |" + $vtab + "9|    |"

$cell = $t.Cell(4,2)
$cell.Range.Text = "25 x 19" + $vtab + "  1    9" + $vtab + "  ----" + $vtab + "2|    |" + $vtab + "5|    |"

$cell = $t.Cell(4,3)
$cell.Range.Text = "19 x 97" + $vtab + "  9    7" + $vtab + "  ----" + $vtab + "1|    |" + $vtab + "9|    |"

$cell = $t.Cell(5,1)
$cell.Range.Text = "94 x 17" + $vtab + "  1    7" + $vtab + "  ----" + $vtab + "9|    |" + $vtab + "4|    |"

$cell = $t.Cell(5,2)
$cell.Range.Text = "63 x 38" + $vtab + "  3    8" + $vtab + "  ----" + $vtab + "6|    |" + $vtab + "3|    |"

$cell = $t.Cell(5,3)
$cell.Range.Text = "98 x 46" + $vtab + "  4    6" + $vtab + "  ----" + $vtab + "9|    |" + $vtab + "8|    |"
